# Apply cryptos.xlsx price/volume update (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.590.82'
$ws.Range('E2').Value = '  +3.95%  '
$ws.Range('D3').Value = '1.743.41'
$ws.Range('E3').Value = '  +4.43%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9994'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.15'
$ws.Range('E5').Value = '  +3.48%  '
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4823'
$ws.Range('E7').Value = '  +1.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2709'
$ws.Range('E8').Value = '  +3.48%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06268'
$ws.Range('E9').Value = '  +1.48%  '
$ws.Range('D10').Value = '1.743.84'
$ws.Range('E10').Value = '  +4.44%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07131'
$ws.Range('E11').Value = '  +1.95%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.88'
$ws.Range('E12').Value = '  +7.05%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6254'
$ws.Range('E13').Value = '  +6.06%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.528'
$ws.Range('E14').Value = '  +3.44%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.55'
$ws.Range('E15').Value = '  +2.86%  '
$ws.Range('E16').Value = '  +0.01%  '
$ws.Range('D17').Value = '26.587.27'
$ws.Range('E17').Value = '  +3.97%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.0000'
$ws.Range('E18').Value = '  +0.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000006900'
$ws.Range('E19').Value = '  +2.26%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.75'
$ws.Range('E20').Value = '  +2.83%  '
$ws.Range('D21').Value = '1.968.32'
$ws.Range('E21').Value = '  +4.33%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.639'
$ws.Range('E22').Value = '  +4.40%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.861'
$ws.Range('E23').Value = '  +0.71%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.383'
$ws.Range('E24').Value = '  +2.29%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '136.33'
$ws.Range('E25').Value = '  -0.32%  '
$ws.Range('E26').Value = '  +2.43%  '
$ws.Range('E27').Value = '  +5.83%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.433'
$ws.Range('E28').Value = '  +3.66%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '106.80'
$ws.Range('E29').Value = '  +1.98%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.009'
$ws.Range('E30').Value = '  +0.27%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.759'
$ws.Range('E31').Value = '  +3.34%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.07893'
$ws.Range('E32').Value = '  +0.29%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04636'
$ws.Range('E33').Value = '  +7.63%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.618'
$ws.Range('E34').Value = '  -0.16%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6447'
$ws.Range('E35').Value = '  +6.54%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9997'
$ws.Range('E36').Value = '  +4.57%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9398'
$ws.Range('E37').Value = '  +2.22%  '
$ws.Range('B38').Value = 'Quant'
$ws.Range('C38').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '113.13'
$ws.Range('E38').Value = '  +15.25%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.997'
$ws.Range('E39').Value = '  +7.81%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.422'
$ws.Range('E40').Value = '  -6.01%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.0000'
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.750'
$ws.Range('E42').Value = '  +17.70%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.01511'
$ws.Range('E43').Value = '  +2.48%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3925'
$ws.Range('E44').Value = '  +4.29%  '
$ws.Range('B45').Value = 'Algorand'
$ws.Range('C45').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1222'
$ws.Range('E45').Value = '  +9.09%  '
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.757'
$ws.Range('E46').Value = '  +8.60%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05334'
$ws.Range('E47').Value = '  +1.27%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.962'
$ws.Range('E48').Value = '  +6.85%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '30.77'
$ws.Range('E49').Value = '  +2.65%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.268'
$ws.Range('E50').Value = '  +5.16%  '
$ws.Range('B51').Value = 'Decentraland'
$ws.Range('C51').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3456'
$ws.Range('E51').Value = '  +3.39%  '
